# Applies a FRED data refresh to the RRPONTSYD workbook:
#  - Appends new daily observations (rows 490-517) to the "Data" sheet
#  - Updates the metadata on the "SeriesInfo" sheet to reflect the new pull

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# New observations to append: (row, date serial, value)
$newRows = @(
    @(490, 45187, 1452.942),
    @(491, 45188, 1453.324),
    @(492, 45189, 1486.984),
    @(493, 45190, 1454.115),
    @(494, 45191, 1427.575),
    @(495, 45194, 1437.31),
    @(496, 45195, 1438.301),
    @(497, 45196, 1442.805),
    @(498, 45197, 1453.366),
    @(499, 45198, 1557.569),
    @(500, 45201, 1365.739),
    @(501, 45202, 1348.465),
    @(502, 45203, 1342.031),
    @(503, 45204, 1265.132),
    @(504, 45205, 1283.461),
    @(505, 45209, 1222.44),
    @(506, 45210, 1239.382),
    @(507, 45211, 1157.319),
    @(508, 45212, 1151.818),
    @(509, 45215, 1108.819),
    @(510, 45216, 1082.502),
    @(511, 45217, 1150.781),
    @(512, 45218, 1114.179),
    @(513, 45219, 1138.756),
    @(514, 45222, 1157.976),
    @(515, 45223, 1097.875),
    @(516, 45224, 1100.617),
    @(517, 45225, 1089.85)
)

$firstNewRow = 490
$lastNewRow = 517
$lastExistingRow = 489

# Stamp the new rows with the same look (date number format, bold border,
# centered alignment) as the rest of the date column before filling values.
$srcRange = $dataSheet.Range("A$lastExistingRow`:B$lastExistingRow")
$destRange = $dataSheet.Range("A$firstNewRow`:B$lastNewRow")
$srcRange.Copy($destRange)

foreach ($row in $newRows) {
    $r = $row[0]
    $dataSheet.Cells.Item($r, 1).Value = $row[1]
    $dataSheet.Cells.Item($r, 2).Value = $row[2]
}

# Update SeriesInfo metadata to reflect the new FRED pull.
# Force these plain "YYYY-MM-DD" cells to stay text (otherwise Excel
# auto-converts the recognizable date string into a date serial number).
$infoSheet.Range("B3").NumberFormat = "@"
$infoSheet.Range("B3").Value = "2023-10-27"
$infoSheet.Range("B4").NumberFormat = "@"
$infoSheet.Range("B4").Value = "2023-10-27"
$infoSheet.Range("B7").NumberFormat = "@"
$infoSheet.Range("B7").Value = "2023-10-26"
$infoSheet.Range("B14").Value = "2023-10-26 13:01:02-05"
$infoSheet.Range("B15").Value = 92
